$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.266.85'
$ws.Range('E2').Value = '  +0.30%  '
$ws.Range('D3').Value = '1.877.26'
$ws.Range('E3').Value = '  -1.46%  '
$ws.Range('E4').Value = '  -0.53%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '246.57'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -2.84%  '
$ws.Range('E6').Value = '  -2.82%  '
$ws.Range('E7').Value = '  -0.59%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '43.65'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +4.42%  '
$ws.Range('E9').Value = '  +0.30%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '53.56'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +2.36%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0739'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -2.53%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0977'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -0.17%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '13.48'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +1.15%  '
$ws.Range('D14').Value = '2.150.50'
$ws.Range('E14').Value = '  -1.42%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.765'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +3.92%  '
$ws.Range('E16').Value = '  -2.24%  '
$ws.Range('D17').Value = '1.879.24'
$ws.Range('E17').Value = '  -1.36%  '
$ws.Range('D18').Value = '35.334.54'
$ws.Range('E18').Value = '  +0.48%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '72.56'
$ws.Range('D19').Style = "Normal"
$ws.Range('D20').Value = '0.0₃0822'
$ws.Range('E20').Value = '  -2.58%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '243.83'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.04%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '12.83'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -1.98%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.98'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -1.44%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.64'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +8.38%  '
$ws.Range('E25').Value = '  -0.49%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.24'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -3.38%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '165.41'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -1.59%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '8.59'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -0.34%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '18.26'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -1.60%  '
$ws.Range('E30').Value = '  -2.25%  '
$ws.Range('E31').Value = '  +1.59%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.67'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +2.05%  '
$ws.Range('E33').Value = '  -1.32%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0591'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -1.15%  '
$ws.Range('E35').Value = '  -2.57%  '
$ws.Range('E36').Value = '  -0.59%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.843'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -1.03%  '
$ws.Range('E38').Value = '  -3.90%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0730'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +10.70%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '17.53'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +1.17%  '
$ws.Range('E41').Value = '  +1.21%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '96.27'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -1.76%  '
$ws.Range('E43').Value = '  -3.49%  '
$ws.Range('D44').Value = '1.306.15'
$ws.Range('E44').Value = '  +0.03%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.38'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -2.13%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0803'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +6.52%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.39'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -1.46%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '11.87'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -4.52%  '
$ws.Range('E50').Value = '  -5.64%  '
$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '42.08'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -2.35%  '
